# Update the expenses/production/revenue "reports" index workbook:
#  - rename Sheet1 -> Reports
#  - widen the Name / generatedAt / filePath columns so the new,
#    longer values are readable
#  - append a row describing the freshly generated expenses report

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Reports"

# Column widths (Office.js/COM ColumnWidth is in points; the OOXML
# <col width> "characters" value this engine emits is ColumnWidth + 5/6,
# so subtract 5/6 here to land exactly on the target character widths).
$ws.Columns.Item(2).ColumnWidth = 30 - 5/6
$ws.Columns.Item(4).ColumnWidth = 26 - 5/6
$ws.Columns.Item(5).ColumnWidth = 30 - 5/6

# New row: id, name, type, generatedAt, filePath
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Expenses Report - 2025-05-05T08-49-07-891Z"
$ws.Range("C2").Value = "expenses"
$ws.Range("D2").Value = "2025-05-05T08:49:07.920Z"
$ws.Range("E2").Value = "/home/runner/workspace/data/reports/expenses_2025-05-05T08-49-07-891Z.xlsx"
